$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.590.96"
$ws.Range("E2").Value = "  -1.09%  "
$ws.Range("D3").Value = "1.672.03"
$ws.Range("E3").Value = "  -2.12%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "'314.38"
$ws.Range("E5").Value = "  -0.20%  "
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").Value = "'0.3959"
$ws.Range("E7").Value = "  -1.56%  "
$ws.Range("D8").Value = "'0.3962"
$ws.Range("E8").Value = "  -2.18%  "
$ws.Range("D9").Value = "'1.004"
$ws.Range("E9").Value = "  +0.09%  "
$ws.Range("D10").Value = "'1.398"
$ws.Range("E10").Value = "  -5.10%  "
$ws.Range("D11").Value = "'50.49"
$ws.Range("E11").Value = "  -6.04%  "
$ws.Range("D12").Value = "'0.08648"
$ws.Range("E12").Value = "  -1.71%  "
$ws.Range("D13").Value = "'25.47"
$ws.Range("E13").Value = "  -3.05%  "
$ws.Range("D14").Value = "'7.313"
$ws.Range("E14").Value = "  -2.60%  "
$ws.Range("D15").Value = "'0.00001316"
$ws.Range("E15").Value = "  -2.04%  "
$ws.Range("D16").Value = "'7.694"
$ws.Range("E16").Value = "  -3.87%  "
$ws.Range("D17").Value = "1.686.47"
$ws.Range("E17").Value = "  +1.42%  "
$ws.Range("D18").Value = "'93.94"
$ws.Range("E18").Value = "  -1.67%  "
$ws.Range("D19").Value = "'0.07021"
$ws.Range("E19").Value = "  -2.13%  "
$ws.Range("D20").Value = "'21.24"
$ws.Range("E20").Value = "  +1.24%  "
$ws.Range("D21").Value = "'7.077"
$ws.Range("E21").Value = "  -2.91%  "
$ws.Range("D22").Value = "'1.002"
$ws.Range("E22").Value = "  -0.23%  "
$ws.Range("D23").Value = "'13.90"
$ws.Range("E23").Value = "  -3.88%  "
$ws.Range("D24").Value = "24.555.21"
$ws.Range("E24").Value = "  -1.31%  "
$ws.Range("D25").Value = "'2.345"
$ws.Range("E25").Value = "  +0.37%  "
$ws.Range("D26").Value = "'2.766"
$ws.Range("E26").Value = "  -4.11%  "
$ws.Range("D27").Value = "'23.01"
$ws.Range("E27").Value = "  -0.20%  "
$ws.Range("D28").Value = "'5.841"
$ws.Range("E28").Value = "  -8.86%  "
$ws.Range("D29").Value = "'159.92"
$ws.Range("E29").Value = "  -1.10%  "
$ws.Range("D30").Value = "'145.82"
$ws.Range("E30").Value = "  +1.22%  "
$ws.Range("D31").Value = "'8.319"
$ws.Range("E31").Value = "  -0.74%  "
$ws.Range("D32").Value = "'2.531"
$ws.Range("E32").Value = "  +10.95%  "
$ws.Range("D33").Value = "1.860.48"
$ws.Range("E33").Value = "  +7.58%  "
$ws.Range("D34").Value = "'0.03083"
$ws.Range("E34").Value = "  -3.18%  "
$ws.Range("D35").Value = "'0.08256"
$ws.Range("E35").Value = "  -5.41%  "
$ws.Range("D36").Value = "'6.912"
$ws.Range("E36").Value = "  -3.77%  "
$ws.Range("E37").Value = "  -2.23%  "
$ws.Range("D38").Value = "'0.9928"
$ws.Range("E38").Value = "  -3.38%  "
$ws.Range("D39").Value = "'0.09651"
$ws.Range("E39").Value = "  +2.07%  "
$ws.Range("D40").Value = "'1.518"
$ws.Range("E40").Value = "  +2.56%  "
$ws.Range("D41").Value = "'10.32"
$ws.Range("E41").Value = "  -4.67%  "
$ws.Range("D42").Value = "'0.7871"
$ws.Range("E42").Value = "  -6.52%  "
$ws.Range("E43").Value = "  -4.53%  "
$ws.Range("D44").Value = "'16.59"
$ws.Range("E44").Value = "  -4.78%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "'0.7097"
$ws.Range("E45").Value = "  -4.54%  "
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").Value = "'2.560"
$ws.Range("E46").Value = "  -5.82%  "
$ws.Range("E47").Value = "  -1.29%  "
$ws.Range("D48").Value = "'0.08640"
$ws.Range("E48").Value = "  +3.04%  "
$ws.Range("E49").Value = "  +0.24%  "
$ws.Range("D50").Value = "'1.328"
$ws.Range("E50").Value = "  -3.50%  "
$ws.Range("D51").Value = "'138.04"
